# Update SMSMY yearly financials worksheet with latest reported figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SMSMY")

# Row 8: Total Revenue
$ws.Range("D8").Value = 4569000
$ws.Range("E8").Value = 3600900
$ws.Range("F8").Value = 3299700
$ws.Range("G8").Value = 4477100
$ws.Range("H8").Value = 4978300
$ws.Range("I8").Value = 5096200
$ws.Range("J8").Value = 6397400

# Row 9: Cost of Revenue
$ws.Range("D9").Value = 3942600
$ws.Range("E9").Value = 3064300
$ws.Range("F9").Value = 2802000
$ws.Range("G9").Value = 3926600
$ws.Range("H9").Value = 4421900
$ws.Range("I9").Value = 4633200
$ws.Range("J9").Value = 5883200

# Row 10: Gross Profit
$ws.Range("D10").Value = 626300
$ws.Range("E10").Value = 536600
$ws.Range("F10").Value = 497700
$ws.Range("G10").Value = 550500
$ws.Range("H10").Value = 556300
$ws.Range("I10").Value = 463000
$ws.Range("J10").Value = 514200

# Row 14: Non Recurring
$ws.Range("D14").Value = -7100
$ws.Range("E14").Value = -2700
$ws.Range("F14").Value = 121800
$ws.Range("H14").Value = 19800
$ws.Range("I14").Value = 215400
$ws.Range("J14").Value = 452700

# Row 15: Others
$ws.Range("D15").Value = 82900
$ws.Range("E15").Value = 79500
$ws.Range("F15").Value = 89400
$ws.Range("G15").Value = 85500
$ws.Range("H15").Value = 83100
$ws.Range("I15").Value = 87400
$ws.Range("J15").Value = 91900

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 4437900
$ws.Range("E17").Value = 3514600
$ws.Range("F17").Value = 3467500
$ws.Range("G17").Value = 4410200
$ws.Range("H17").Value = 4929600
$ws.Range("I17").Value = 5441500
$ws.Range("J17").Value = 6869800

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = 131100
$ws.Range("E18").Value = 86200
$ws.Range("F18").Value = -167700
$ws.Range("G18").Value = 66900
$ws.Range("H18").Value = 48700
$ws.Range("I18").Value = -345300
$ws.Range("J18").Value = -472400

# Row 20: Total Other Income/Expenses Net
$ws.Range("D20").Value = 61800
$ws.Range("E20").Value = 51100
$ws.Range("F20").Value = 12000
$ws.Range("G20").Value = 36400
$ws.Range("H20").Value = 6400
$ws.Range("I20").Value = 12200
$ws.Range("J20").Value = 28700

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 276000
$ws.Range("E21").Value = 216900
$ws.Range("F21").Value = -66200
$ws.Range("G21").Value = 189100
$ws.Range("H21").Value = 142800
$ws.Range("I21").Value = -245500
$ws.Range("J21").Value = -351700

# Row 22: Interest Expense
$ws.Range("D22").Value = 2100
$ws.Range("F22").Value = 3600
$ws.Range("G22").Value = 6400
$ws.Range("H22").Value = 10700
$ws.Range("I22").Value = 12700
$ws.Range("J22").Value = 12800

# Row 23: Income Before Tax
$ws.Range("D23").Value = 190800
$ws.Range("E23").Value = 135100
$ws.Range("F23").Value = -159300
$ws.Range("G23").Value = 96900
$ws.Range("H23").Value = 44400
$ws.Range("I23").Value = -345700
$ws.Range("J23").Value = -456500

# Row 24: Income Tax Expense
$ws.Range("D24").Value = 46800
$ws.Range("E24").Value = -8900
$ws.Range("F24").Value = -6200
$ws.Range("G24").Value = 19200
$ws.Range("H24").Value = 32800
$ws.Range("I24").Value = -15100
$ws.Range("J24").Value = -15800

# Row 26: Income After Tax
$ws.Range("D26").Value = 144000
$ws.Range("E26").Value = 144000
$ws.Range("F26").Value = -153200
$ws.Range("G26").Value = 77700
$ws.Range("H26").Value = 11500
$ws.Range("I26").Value = -330600
$ws.Range("J26").Value = -440800

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = 144000
$ws.Range("E27").Value = 144000
$ws.Range("F27").Value = -153200
$ws.Range("G27").Value = 77700
$ws.Range("H27").Value = 11500
$ws.Range("I27").Value = -330600
$ws.Range("J27").Value = -440800

# Row 29: Discontinued Operations
$ws.Range("H29").Value = -74400

# Row 32: Other Items
$ws.Range("D32").Value = -61800
$ws.Range("E32").Value = -51100
$ws.Range("F32").Value = -12000
$ws.Range("G32").Value = -36400
$ws.Range("H32").Value = -6400
$ws.Range("I32").Value = -12200
$ws.Range("J32").Value = -28700

# Row 33: Net Income
$ws.Range("D33").Value = 144000
$ws.Range("E33").Value = 144000
$ws.Range("F33").Value = -153200
$ws.Range("G33").Value = 77800
$ws.Range("H33").Value = -62900
$ws.Range("I33").Value = -330600
$ws.Range("J33").Value = -440800

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = 144000
$ws.Range("E35").Value = 144000
$ws.Range("F35").Value = -153200
$ws.Range("G35").Value = 77800
$ws.Range("H35").Value = -62900
$ws.Range("I35").Value = -330600
$ws.Range("J35").Value = -440800

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 239900
$ws.Range("E41").Value = 267800
$ws.Range("F41").Value = 175700
$ws.Range("G41").Value = 223600
$ws.Range("H41").Value = 40500
$ws.Range("I41").Value = 33200
$ws.Range("J41").Value = 36400

# Row 42: Short Term Investments
$ws.Range("D42").Value = 8000
$ws.Range("E42").Value = 16300
$ws.Range("F42").Value = 6300
$ws.Range("G42").Value = 6400
$ws.Range("H42").Value = 5300
$ws.Range("I42").Value = 25600
$ws.Range("J42").Value = 20800

# Row 43: Net Receivables
$ws.Range("D43").Value = 310700
$ws.Range("E43").Value = 584300
$ws.Range("F43").Value = 264700
$ws.Range("G43").Value = 263600
$ws.Range("H43").Value = 325200
$ws.Range("I43").Value = 640900
$ws.Range("J43").Value = 349800

# Row 44: Inventory
$ws.Range("D44").Value = 401200
$ws.Range("E44").Value = 569800
$ws.Range("F44").Value = 281800
$ws.Range("G44").Value = 337100
$ws.Range("H44").Value = 385600
$ws.Range("I44").Value = 822900
$ws.Range("J44").Value = 586900

# Row 45: Other Current Assets
$ws.Range("D45").Value = 22900
$ws.Range("E45").Value = 24100
$ws.Range("F45").Value = 36200
$ws.Range("G45").Value = 22800
$ws.Range("H45").Value = 26600
$ws.Range("I45").Value = 56700
$ws.Range("J45").Value = 17300

# Row 46: Total Current Assets
$ws.Range("D46").Value = 982600
$ws.Range("E46").Value = 868600
$ws.Range("F46").Value = 764600
$ws.Range("G46").Value = 853400
$ws.Range("H46").Value = 783200
$ws.Range("I46").Value = 813100
$ws.Range("J46").Value = 1011100

# Row 47: Long Term Investments
$ws.Range("D47").Value = 201000
$ws.Range("E47").Value = 159800
$ws.Range("F47").Value = 146000
$ws.Range("G47").Value = 220600
$ws.Range("H47").Value = 227900
$ws.Range("I47").Value = 307300
$ws.Range("J47").Value = 283400

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 817700
$ws.Range("E48").Value = 1373800
$ws.Range("F48").Value = 697000
$ws.Range("G48").Value = 730000
$ws.Range("H48").Value = 638900
$ws.Range("I48").Value = 904000
$ws.Range("J48").Value = 691300

# Row 49: Goodwill
$ws.Range("D49").Value = 140500
$ws.Range("E49").Value = 224100
$ws.Range("F49").Value = 120100
$ws.Range("G49").Value = 162200
$ws.Range("H49").Value = 152700
$ws.Range("I49").Value = 560100
$ws.Range("J49").Value = 405300

# Row 52: Other Assets
$ws.Range("D52").Value = 123500
$ws.Range("E52").Value = 121100
$ws.Range("F52").Value = 91300
$ws.Range("G52").Value = 72700
$ws.Range("H52").Value = 71800
$ws.Range("I52").Value = 88000
$ws.Range("J52").Value = 91600

# Row 54: Total Assets
$ws.Range("D54").Value = 2265300
$ws.Range("E54").Value = 1940700
$ws.Range("F54").Value = 1818900
$ws.Range("G54").Value = 2038900
$ws.Range("H54").Value = 1874500
$ws.Range("I54").Value = 2064100
$ws.Range("J54").Value = 2482600

# Row 57: Accounts Payable
$ws.Range("D57").Value = 268400
$ws.Range("E57").Value = 514300
$ws.Range("F57").Value = 183600
$ws.Range("G57").Value = 225600
$ws.Range("H57").Value = 267600
$ws.Range("I57").Value = 710600
$ws.Range("J57").Value = 343300

# Row 58: Short/Current Long Term Debt
$ws.Range("D58").Value = 1100
$ws.Range("I58").Value = 16400
$ws.Range("J58").Value = 9700

# Row 59: Other Current Liabilities
$ws.Range("D59").Value = 292500
$ws.Range("E59").Value = 230500
$ws.Range("F59").Value = 188700
$ws.Range("G59").Value = 202800
$ws.Range("H59").Value = 211400
$ws.Range("I59").Value = 170200
$ws.Range("J59").Value = 166500

# Row 60: Total Current Liabilities
$ws.Range("D60").Value = 562000
$ws.Range("E60").Value = 420000
$ws.Range("F60").Value = 373800
$ws.Range("G60").Value = 428700
$ws.Range("H60").Value = 479400
$ws.Range("I60").Value = 475000
$ws.Range("J60").Value = 519400

# Row 61: Long Term Debt
$ws.Range("D61").Value = 27900
$ws.Range("E61").Value = 2400
$ws.Range("F61").Value = 2900
$ws.Range("G61").Value = 1100
$ws.Range("H61").Value = 10200
$ws.Range("I61").Value = 133800
$ws.Range("J61").Value = 233400

# Row 62: Other Liabilities
$ws.Range("D62").Value = 126900
$ws.Range("E62").Value = 126200
$ws.Range("F62").Value = 145700
$ws.Range("G62").Value = 114200
$ws.Range("H62").Value = 87400
$ws.Range("I62").Value = 119000
$ws.Range("J62").Value = 114000

# Row 66: Total Liabilities
$ws.Range("D66").Value = 716800
$ws.Range("E66").Value = 548600
$ws.Range("F66").Value = 522400
$ws.Range("G66").Value = 544100
$ws.Range("H66").Value = 577000
$ws.Range("I66").Value = 699200
$ws.Range("J66").Value = 866900

# Row 72: Retained Earnings
$ws.Range("D72").Value = -369200
$ws.Range("E72").Value = -460200
$ws.Range("F72").Value = -585200
$ws.Range("G72").Value = -403200
$ws.Range("H72").Value = -448600
$ws.Range("I72").Value = -533200
$ws.Range("J72").Value = -64700

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 1548500
$ws.Range("E76").Value = 1392100
$ws.Range("F76").Value = 1296500
$ws.Range("G76").Value = 1494800
$ws.Range("H76").Value = 1297500
$ws.Range("I76").Value = 1364900
$ws.Range("J76").Value = 1615700

# Row 81: Net Income
$ws.Range("D81").Value = 144000
$ws.Range("E81").Value = 144000
$ws.Range("F81").Value = -153200
$ws.Range("G81").Value = 77800
$ws.Range("H81").Value = -62900
$ws.Range("I81").Value = -330600
$ws.Range("J81").Value = -440800

# Row 83: Depreciation
$ws.Range("D83").Value = 82900
$ws.Range("E83").Value = 79500
$ws.Range("F83").Value = 89400
$ws.Range("G83").Value = 85700
$ws.Range("H83").Value = 87700
$ws.Range("I83").Value = 87400
$ws.Range("J83").Value = 91900

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 178400
$ws.Range("E89").Value = 188500
$ws.Range("F89").Value = 92900
$ws.Range("G89").Value = 210900
$ws.Range("H89").Value = 148600
$ws.Range("I89").Value = 210300
$ws.Range("J89").Value = 204900

# Row 91: Capital Expenditures
$ws.Range("D91").Value = -124600
$ws.Range("E91").Value = -89500
$ws.Range("F91").Value = -77000
$ws.Range("G91").Value = -67400
$ws.Range("H91").Value = -45400
$ws.Range("I91").Value = -105400
$ws.Range("J91").Value = -114000

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -185600
$ws.Range("E94").Value = -44600
$ws.Range("F94").Value = -67800
$ws.Range("G94").Value = 29500
$ws.Range("H94").Value = -600
$ws.Range("I94").Value = -84100
$ws.Range("J94").Value = -242400

# Row 96: Dividends Paid
$ws.Range("D96").Value = -75600
$ws.Range("E96").Value = -44700
$ws.Range("F96").Value = -33100
$ws.Range("G96").Value = -37600
$ws.Range("I96").Value = -14400
$ws.Range("J96").Value = -49000

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -25900
$ws.Range("E100").Value = -47800
$ws.Range("F100").Value = -78900
$ws.Range("G100").Value = -71000
$ws.Range("H100").Value = -141300
$ws.Range("I100").Value = -131900
$ws.Range("J100").Value = -46300

# Row 101: Effect Of Exchange Rate Changes 
$ws.Range("D101").Value = 5300
$ws.Range("E101").Value = -4000
$ws.Range("F101").Value = 5900
$ws.Range("G101").Value = 13700
$ws.Range("J101").Value = 3000

# Row 102: Change In Cash and Cash Equivalents 
$ws.Range("D102").Value = -27900
$ws.Range("E102").Value = 92100
$ws.Range("F102").Value = -47900
$ws.Range("G102").Value = 183100
$ws.Range("H102").Value = 7300
$ws.Range("I102").Value = -3200
$ws.Range("J102").Value = -80700
